# Update "want to go" counts (column F) on several sheets to reflect
# newly generated output (gh-pages data refresh at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 123
$wsExhibit.Range("F3").Value = 431

# Sheet "演出" (Performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 69
$wsShow.Range("F3").Value = 28

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 123
$wsAll.Range("F3").Value = 69
$wsAll.Range("F4").Value = 431
$wsAll.Range("F8").Value = 28
